$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the data set. It lands at row 445,
# pushing the existing rows 445-457 down to 446-458.
$ws.Rows.Item(445).Insert()

# Populate the newly inserted row 445 with the new record's values.
$ws.Cells.Item(445, 1).Value = 6
$ws.Cells.Item(445, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(445, 3).Value = "Metropolitana"
$ws.Cells.Item(445, 4).Value = 44890
$ws.Cells.Item(445, 5).Value = 13
$ws.Cells.Item(445, 6).Value = "Fruta"
$ws.Cells.Item(445, 7).Value = 100107
$ws.Cells.Item(445, 8).Value = "Otros"
$ws.Cells.Item(445, 9).Value = 100107011
$ws.Cells.Item(445, 10).Value = "Tuna"
$ws.Cells.Item(445, 11).Value = "Sin especificar"
$ws.Cells.Item(445, 12).Value = "Especial"
$ws.Cells.Item(445, 13).Value = 60
$ws.Cells.Item(445, 14).Value = 35000
$ws.Cells.Item(445, 15).Value = 35000
$ws.Cells.Item(445, 16).Value = 35000
$ws.Cells.Item(445, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(445, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(445, 19).Value = 1944
$ws.Cells.Item(445, 20).Value = 18
